$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A137").Value = "MRN:  JH16121935"
$ws.Range("B137").Value = 0.1685060236161159
$ws.Range("C137").Value = 0.3170623490466765
$ws.Range("D137").Value = 0.2257615077889461
$ws.Range("E137").Value = 0.2298821502492125
$ws.Range("F137").Value = 0.0962215813268736
$ws.Range("G137").Value = 0.2484848305735292

$ws.Range("A138").Value = "MRN:  JH16121937"
$ws.Range("B138").Value = 0.1890888150049136
$ws.Range("C138").Value = 0.1388648298213223
$ws.Range("D138").Value = 0.2331879231993602
$ws.Range("E138").Value = 0.08975929591124644
$ws.Range("F138").Value = 0.2067592471613154
$ws.Range("G138").Value = 0.2331879231993602
